$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-18 from 2023-09-01 to 2023-09-05
$ws.Range("C2:C18").Value = 45174
